# Change the D1:F1 header cells from numeric dates (stored with a custom
# "YYYY-MM-DD" number format) into plain text strings showing the dates
# in "YYYY.MM.DD" form. The cells keep looking like the rest of the header
# row (bold font, thin border, centered/top aligned) - they simply no
# longer need the custom date number format, so that style becomes
# unused.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @{ "D1" = "2025.12.01"; "E1" = "2025.12.08"; "F1" = "2025.12.15" }

foreach ($addr in @("D1", "E1", "F1")) {
    $cell = $ws.Range($addr)

    # Mark the cell as Text first so assigning the value afterwards is
    # stored verbatim as a string instead of being re-parsed back into a
    # date serial number.
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$addr]

    # Drop back to the plain (non-date) look shared by the rest of the
    # header row: bold font, thin border on all sides, centered
    # horizontally and top-aligned vertically - same as A1:C1.
    $cell.ClearFormats()
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108  # xlCenter
    $cell.VerticalAlignment = -4160    # xlTop
    $cell.Borders.LineStyle = 1        # xlContinuous
}
